$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert "minimum number of steps required to sort an array":
# remove row 22 entirely, and turn A21 back into the plain numeric value -20
# (was a shared-string entry referencing "( 20)quick  sort ...").

$ws.Rows.Item(22).Delete()
$ws.Range("A21").Value = -20

$ws.Range("A21").Select()
